$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value pairs derived from the crypto price-ticker refresh.
# NumberFormat="@" + ClearFormats keeps numeric-looking strings (e.g. "576.17")
# stored as text, matching the source data which is never a true number.
$updates = @(
    @{Cell='D2'; Value='64.000.49'}
    @{Cell='E2'; Value='  +0.18%  '}
    @{Cell='D3'; Value='2.760.42'}
    @{Cell='E3'; Value='  -0.33%  '}
    @{Cell='E4'; Value='  +0.07%  '}
    @{Cell='D5'; Value='576.17'}
    @{Cell='E5'; Value='  -1.50%  '}
    @{Cell='D6'; Value='158.98'}
    @{Cell='E6'; Value='  -1.80%  '}
    @{Cell='D7'; Value='0.998'}
    @{Cell='E7'; Value='  +0.21%  '}
    @{Cell='E8'; Value='  -3.36%  '}
    @{Cell='E9'; Value='  -3.67%  '}
    @{Cell='E11'; Value='  -14.72%  '}
    @{Cell='E12'; Value='  -3.25%  '}
    @{Cell='D13'; Value='3.248.92'}
    @{Cell='E13'; Value='  -0.43%  '}
    @{Cell='D14'; Value='26.87'}
    @{Cell='E14'; Value='  -2.45%  '}
    @{Cell='D15'; Value='63.699.50'}
    @{Cell='E15'; Value='  -0.21%  '}
    @{Cell='E16'; Value='  -4.92%  '}
    @{Cell='D17'; Value='2.763.11'}
    @{Cell='E17'; Value='  -0.81%  '}
    @{Cell='D18'; Value='12.11'}
    @{Cell='E18'; Value='  -1.75%  '}
    @{Cell='D19'; Value='4.83'}
    @{Cell='E19'; Value='  -3.51%  '}
    @{Cell='D20'; Value='355.59'}
    @{Cell='E20'; Value='  -3.39%  '}
    @{Cell='D21'; Value='6.65'}
    @{Cell='E21'; Value='  -6.08%  '}
    @{Cell='D22'; Value='0.998'}
    @{Cell='E22'; Value='  -0.42%  '}
    @{Cell='E23'; Value='  -5.47%  '}
    @{Cell='D24'; Value='65.00'}
    @{Cell='E24'; Value='  -3.42%  '}
    @{Cell='E25'; Value='  -3.65%  '}
    @{Cell='D26'; Value='8.53'}
    @{Cell='E26'; Value='  -2.04%  '}
    @{Cell='E27'; Value='  +0.39%  '}
    @{Cell='E28'; Value='  -5.92%  '}
    @{Cell='D29'; Value='7.36'}
    @{Cell='E29'; Value='  +0.09%  '}
    @{Cell='D30'; Value='1.95'}
    @{Cell='E30'; Value='  -4.12%  '}
    @{Cell='E31'; Value='  -0.56%  '}
    @{Cell='D32'; Value='170.07'}
    @{Cell='E32'; Value='  -1.66%  '}
    @{Cell='E33'; Value='  -3.06%  '}
    @{Cell='D34'; Value='20.13'}
    @{Cell='E34'; Value='  -3.31%  '}
    @{Cell='E35'; Value='  +0.07%  '}
    @{Cell='E36'; Value='  -1.38%  '}
    @{Cell='D37'; Value='1.79'}
    @{Cell='E37'; Value='  -2.57%  '}
    @{Cell='D38'; Value='0.999'}
    @{Cell='E38'; Value='  -3.92%  '}
    @{Cell='D39'; Value='350.01'}
    @{Cell='E39'; Value='  +1.69%  '}
    @{Cell='E40'; Value='  -0.54%  '}
    @{Cell='D41'; Value='4.16'}
    @{Cell='E41'; Value='  -2.67%  '}
    @{Cell='E42'; Value='  -2.25%  '}
    @{Cell='D43'; Value='21.41'}
    @{Cell='E43'; Value='  -4.71%  '}
    @{Cell='D44'; Value='21.76'}
    @{Cell='E44'; Value='  -5.08%  '}
    @{Cell='D45'; Value='0.0585'}
    @{Cell='D46'; Value='137.83'}
    @{Cell='E46'; Value='  -0.66%  '}
    @{Cell='E47'; Value='  -3.27%  '}
    @{Cell='D48'; Value='0.0254'}
    @{Cell='E48'; Value='  -3.17%  '}
    @{Cell='E49'; Value='  -1.76%  '}
    @{Cell='E50'; Value='  +0.27%  '}
    @{Cell='E51'; Value='  +0.19%  '}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = '@'
    $rng.Value = $u.Value
    $rng.ClearFormats()
}
